# Rewrite the CloudFront-hosted image URLs baked into each picture's
# description (alt text) to the local "images/..." paths used by the
# published site (gh-pages build).
#
# Mapping (old descr -> new descr):
#   https://deidt7p41jzcy.cloudfront.net/jejunostomy_qrcode.png  -> images/nutrition_jejunostomy_qrcode.png
#   https://deidt7p41jzcy.cloudfront.net/Eso_Anatomy_Labels.png  -> images/Eso_Anatomy_Labels.png
#   https://deidt7p41jzcy.cloudfront.net/protein_shakes.png      -> images/nutrition_protein_shakes.png
#   https://deidt7p41jzcy.cloudfront.net/jtube_ai.png            -> images/nutrition_jtube.png
#   https://deidt7p41jzcy.cloudfront.net/gtube_ai.png            -> images/nutrition_gtube.png

$p = $ppt.ActivePresentation

function Set-PictureDescr {
    param($Slide, $OldDescr, $NewDescr)

    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $shp = $Slide.Shapes.Item($i)
        if ($shp.Type -eq 13 -and $shp.AlternativeText -eq $OldDescr) {
            $shp.AlternativeText = $NewDescr
            return $true
        }
    }
    return $false
}

# --- Slide 2: Eso_Anatomy_Labels.png ---------------------------------------
$s2 = $p.Slides.Item(2)
Set-PictureDescr $s2 "https://deidt7p41jzcy.cloudfront.net/Eso_Anatomy_Labels.png" "images/Eso_Anatomy_Labels.png" | Out-Null

# --- Slide 3: protein_shakes.png --------------------------------------------
$s3 = $p.Slides.Item(3)
Set-PictureDescr $s3 "https://deidt7p41jzcy.cloudfront.net/protein_shakes.png" "images/nutrition_protein_shakes.png" | Out-Null

# --- Slide 4: jtube_ai.png AND gtube_ai.png ---------------------------------
# Both pictures on this slide share the same underlying shape id (id="0" /
# name="Picture 1" in the source XML), so they cannot be told apart through
# ordinary Shapes.Item(i).AlternativeText get/set calls (those always resolve
# to the physically-first picture). Work around the collision by duplicating
# each picture in turn -- a pasted shape is minted with a fresh, unique id and
# becomes independently addressable -- then remove the original ambiguous
# shapes and apply the correct alt text to each duplicate.
$s4 = $p.Slides.Item(4)

$jtubeOld = $s4.Shapes.Item(2)
$jtubeOld.Copy()
$s4.Shapes.Paste() | Out-Null
$jtubeNew = $s4.Shapes.Item($s4.Shapes.Count)
$s4.Shapes.Item(2).Delete()

$gtubeOld = $s4.Shapes.Item(2)
$gtubeOld.Copy()
$s4.Shapes.Paste() | Out-Null
$gtubeNew = $s4.Shapes.Item($s4.Shapes.Count)
$s4.Shapes.Item(2).Delete()

$jtubeNew.AlternativeText = "images/nutrition_jtube.png"
$gtubeNew.AlternativeText = "images/nutrition_gtube.png"

# --- Slide 8: jtube_ai.png ---------------------------------------------------
$s8 = $p.Slides.Item(8)
Set-PictureDescr $s8 "https://deidt7p41jzcy.cloudfront.net/jtube_ai.png" "images/nutrition_jtube.png" | Out-Null

# --- Slide 11: jejunostomy_qrcode.png ---------------------------------------
$s11 = $p.Slides.Item(11)
Set-PictureDescr $s11 "https://deidt7p41jzcy.cloudfront.net/jejunostomy_qrcode.png" "images/nutrition_jejunostomy_qrcode.png" | Out-Null
